{"js": "// Update the worksheet date and all equation cells to the new day's values.\nconst replacements = [\n  [\"2024-08-04 Sunday\", \"2024-08-05 Monday\"],\n  [\"865\u00d76=\", \"939\u00d77=\"],\n  [\"253\u00d79=\", \"679\u00d73=\"],\n  [\"827\u00d79=\", \"484\u00d79=\"],\n  [\"391\u00d75=\", \"539\u00d78=\"],\n  [\"742\u00d77=\", \"437\u00d77=\"],\n  [\"533\u00d73=\", \"690\u00d72=\"],\n  [\"804\u00d78=\", \"593\u00d76=\"],\n  [\"946\u00d77=\", \"556\u00d78=\"],\n  [\"357\u00d72=\", \"651\u00d77=\"],\n  [\"931\u00d75=\", \"215\u00d76=\"],\n  [\"328\u00d78=\", \"170\u00d73=\"],\n  [\"398\u00d74=\", \"199\u00d76=\"],\n  [\"845\u00d78=\", \"294\u00d73=\"],\n  [\"730\u00d72=\", \"506\u00d79=\"],\n  [\"251\u00d78=\", \"490\u00d74=\"],\n  [\"890\u00d78=\", \"231\u00d73=\"],\n  [\"357\u00d73=\", \"359\u00d72=\"],\n  [\"204\u00d79=\", \"770\u00d73=\"],\n  [\"746\u00d72=\", \"531\u00d79=\"],\n  [\"590\u00d75=\", \"347\u00d74=\"],\n  [\"813\u00d73=\", \"264\u00d73=\"],\n  [\"634\u00d76=\", \"311\u00d75=\"],\n  [\"508\u00d77=\", \"558\u00d77=\"],\n  [\"626\u00d76=\", \"589\u00d77=\"],\n  [\"469\u00d74=\", \"951\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all equation cells to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"2024-08-04 Sunday\"; new = \"2024-08-05 Monday\"},\n    @{old = \"865\u00d76=\"; new = \"939\u00d77=\"},\n    @{old = \"253\u00d79=\"; new = \"679\u00d73=\"},\n    @{old = \"827\u00d79=\"; new = \"484\u00d79=\"},\n    @{old = \"391\u00d75=\"; new = \"539\u00d78=\"},\n    @{old = \"742\u00d77=\"; new = \"437\u00d77=\"},\n    @{old = \"533\u00d73=\"; new = \"690\u00d72=\"},\n    @{old = \"804\u00d78=\"; new = \"593\u00d76=\"},\n    @{old = \"946\u00d77=\"; new = \"556\u00d78=\"},\n    @{old = \"357\u00d72=\"; new = \"651\u00d77=\"},\n    @{old = \"931\u00d75=\"; new = \"215\u00d76=\"},\n    @{old = \"328\u00d78=\"; new = \"170\u00d73=\"},\n    @{old = \"398\u00d74=\"; new = \"199\u00d76=\"},\n    @{old = \"845\u00d78=\"; new = \"294\u00d73=\"},\n    @{old = \"730\u00d72=\"; new = \"506\u00d79=\"},\n    @{old = \"251\u00d78=\"; new = \"490\u00d74=\"},\n    @{old = \"890\u00d78=\"; new = \"231\u00d73=\"},\n    @{old = \"357\u00d73=\"; new = \"359\u00d72=\"},\n    @{old = \"204\u00d79=\"; new = \"770\u00d73=\"},\n    @{old = \"746\u00d72=\"; new = \"531\u00d79=\"},\n    @{old = \"590\u00d75=\"; new = \"347\u00d74=\"},\n    @{old = \"813\u00d73=\"; new = \"264\u00d73=\"},\n    @{old = \"634\u00d76=\"; new = \"311\u00d75=\"},\n    @{old = \"508\u00d77=\"; new = \"558\u00d77=\"},\n    @{old = \"626\u00d76=\"; new = \"589\u00d77=\"},\n    @{old = \"469\u00d74=\"; new = \"951\u00d73=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.new\n    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
